$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, DAMSLTag (I), DialogAct (J) updates
$updates = @(
    ,@(2, 'b', 'Acknowledge (Backchannel)')
    ,@(11, 'aa', 'Agree/Accept')
    ,@(15, 'b', 'Acknowledge (Backchannel)')
    ,@(31, 'aa', 'Agree/Accept')
    ,@(46, 'aa', 'Agree/Accept')
    ,@(51, 'sv', 'Statement-opinion')
    ,@(62, 'aa', 'Agree/Accept')
    ,@(65, 'b', 'Acknowledge (Backchannel)')
    ,@(70, 'sv', 'Statement-opinion')
    ,@(71, 'sv', 'Statement-opinion')
    ,@(77, '%', 'Uninterpretable')
    ,@(98, 'sv', 'Statement-opinion')
    ,@(101, 'aa', 'Agree/Accept')
    ,@(108, 'aa', 'Agree/Accept')
    ,@(109, 'b', 'Acknowledge (Backchannel)')
    ,@(110, 'sv', 'Statement-opinion')
    ,@(112, '%', 'Uninterpretable')
    ,@(119, 'aa', 'Agree/Accept')
    ,@(124, 'sv', 'Statement-opinion')
    ,@(132, 'sv', 'Statement-opinion')
    ,@(137, 'sd', 'Statement-non-opinion')
    ,@(141, 'aa', 'Agree/Accept')
    ,@(145, 'b', 'Acknowledge (Backchannel)')
    ,@(152, 'sd', 'Statement-non-opinion')
    ,@(153, 'sv', 'Statement-opinion')
    ,@(155, 'sv', 'Statement-opinion')
    ,@(157, 'b', 'Acknowledge (Backchannel)')
    ,@(161, 'sv', 'Statement-opinion')
    ,@(162, 'b', 'Acknowledge (Backchannel)')
    ,@(165, '%', 'Uninterpretable')
    ,@(180, 'sd', 'Statement-non-opinion')
    ,@(195, 'aa', 'Agree/Accept')
    ,@(196, 'sd', 'Statement-non-opinion')
    ,@(199, 'aa', 'Agree/Accept')
    ,@(205, 'sd', 'Statement-non-opinion')
    ,@(214, 'b', 'Acknowledge (Backchannel)')
    ,@(223, 'aa', 'Agree/Accept')
    ,@(226, 'aa', 'Agree/Accept')
    ,@(228, 'aa', 'Agree/Accept')
    ,@(232, 'ba', 'Appreciation')
    ,@(241, 'ba', 'Appreciation')
    ,@(242, 'sd', 'Statement-non-opinion')
    ,@(243, 'aa', 'Agree/Accept')
    ,@(254, 'sd', 'Statement-non-opinion')
    ,@(264, 'sd', 'Statement-non-opinion')
    ,@(265, 'sd', 'Statement-non-opinion')
    ,@(270, 'sd', 'Statement-non-opinion')
    ,@(284, 'sd', 'Statement-non-opinion')
    ,@(307, 'sd', 'Statement-non-opinion')
    ,@(316, 'sv', 'Statement-opinion')
    ,@(317, 'b', 'Acknowledge (Backchannel)')
    ,@(318, 'ba', 'Appreciation')
    ,@(319, 'sv', 'Statement-opinion')
    ,@(320, 'sv', 'Statement-opinion')
    ,@(321, 'ba', 'Appreciation')
    ,@(324, 'sv', 'Statement-opinion')
    ,@(325, 'aa', 'Agree/Accept')
    ,@(330, 'sd', 'Statement-non-opinion')
    ,@(333, 'ba', 'Appreciation')
    ,@(340, 'sd', 'Statement-non-opinion')
    ,@(354, 'b', 'Acknowledge (Backchannel)')
    ,@(370, 'sd', 'Statement-non-opinion')
    ,@(371, 'ba', 'Appreciation')
    ,@(372, 'sd', 'Statement-non-opinion')
    ,@(375, 'sv', 'Statement-opinion')
    ,@(380, 'aa', 'Agree/Accept')
    ,@(392, 'sv', 'Statement-opinion')
    ,@(395, 'ba', 'Appreciation')
    ,@(405, 'sd', 'Statement-non-opinion')
    ,@(417, 'sv', 'Statement-opinion')
    ,@(418, 'sv', 'Statement-opinion')
    ,@(419, 'sd', 'Statement-non-opinion')
    ,@(421, 'sd', 'Statement-non-opinion')
    ,@(433, 'sd', 'Statement-non-opinion')
    ,@(436, 'sd', 'Statement-non-opinion')
    ,@(437, 'sv', 'Statement-opinion')
    ,@(441, 'sv', 'Statement-opinion')
    ,@(444, 'sd', 'Statement-non-opinion')
    ,@(451, 'sv', 'Statement-opinion')
    ,@(463, 'sv', 'Statement-opinion')
    ,@(467, '%', 'Uninterpretable')
    ,@(477, 'sv', 'Statement-opinion')
    ,@(491, 'aa', 'Agree/Accept')
    ,@(502, 'aa', 'Agree/Accept')
    ,@(507, 'sd', 'Statement-non-opinion')
    ,@(519, 'aa', 'Agree/Accept')
    ,@(531, 'sd', 'Statement-non-opinion')
    ,@(558, 'sv', 'Statement-opinion')
    ,@(563, 'aa', 'Agree/Accept')
    ,@(572, 'sd', 'Statement-non-opinion')
    ,@(577, 'sd', 'Statement-non-opinion')
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 9).Value = $u[1]
    $ws.Cells.Item($r, 10).Value = $u[2]
}
